$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp Mapping")
$ws.Select()
